# The commit trims each of the "amount" values in column B down to their
# first three digits (1234 -> 123, 4321 -> 432, 8765 -> 876, 5434 -> 543,
# 1479 -> 147, 3092 -> 309, 2947 -> 294). Rows 5, 7 and 10 are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 123
$ws.Range("B2").Value = 432
$ws.Range("B3").Value = 876
$ws.Range("B4").Value = 543
$ws.Range("B6").Value = 147
$ws.Range("B8").Value = 309
$ws.Range("B9").Value = 294
